$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.625.24"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").Value = "2.556.07"
$ws.Range("E3").Value = "  +4.19%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.34"
$ws.Range("E5").Value = "  +1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.33"
$ws.Range("E6").Value = "  -4.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  -6.15%  "

$ws.Range("D9").Value = "2.580.40"
$ws.Range("E9").Value = "  +3.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.56"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("E13").Value = "  +0.51%  "

$ws.Range("D14").Value = "3.027.97"
$ws.Range("E14").Value = "  +5.87%  "

$ws.Range("D15").Value = "59.833.68"
$ws.Range("E15").Value = "  +3.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.40"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "2.601.39"
$ws.Range("E18").Value = "  +6.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.69"
$ws.Range("E20").Value = "  +4.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.61"
$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.416"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.705.26"
$ws.Range("E26").Value = "  +6.92%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").Value = "0.0₃0838"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.23"
$ws.Range("E32").Value = "  +2.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.09"
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.66"
$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.94"
$ws.Range("E36").Value = "  +3.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  +22.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").Value = "  +3.65%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.43"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "294.27"
$ws.Range("E43").Value = "  +4.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0560"
$ws.Range("E44").Value = "  +3.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.616"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0992"
$ws.Range("E46").Value = "  -1.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.59"
$ws.Range("E48").Value = "  +7.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("D50").Value = "2.027.58"
$ws.Range("E50").Value = "  +6.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0231"
$ws.Range("E51").Value = "  -0.08%  "
